$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "latest update" timestamp, inserted at the top of the D-column history,
# shifting the previously existing values down (D2:D15 <- new value,
# D16:D29 <- old D2:D15 value, D30:D43 <- old D16:D29 value).
$newTimestamp = 44303.53956131134
$shift1 = 44303.5181042824
$shift2 = 44303.49667890046

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newTimestamp
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $shift1
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $shift2
}
